$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.441.00"
$ws.Range("D3").Value = "3.575.14"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.65%  "
$ws.Range("D7").Value = "3.572.34"
$ws.Range("E7").Value = "  +2.07%  "
$ws.Range("E8").Value = "  +1.80%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.215"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.650"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000311"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").Value = "4.140.55"
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("D16").Value = "70.534.33"
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("E17").Value = "  +4.30%  "
$ws.Range("D18").Value = "3.569.64"
$ws.Range("E18").Value = "  +1.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "575.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "94.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("E27").Value = "  -1.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.01%  "
$ws.Range("E29").Value = "  +4.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("E33").Value = "  +1.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "64.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.74"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +20.50%  "
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.26%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "3.808.38"
$ws.Range("E37").Value = "  +13.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.407"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "522.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.25%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "0.0₃0788"
$ws.Range("E42").Value = "  +3.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.62%  "
$ws.Range("E44").Value = "  +3.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0458"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.13%  "
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.140"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.68%  "
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("E51").Value = "  +7.07%  "
